$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2: existing row, Target cluster changes from "M2" to "ECs", plus updated metrics
$ws.Range("A2").Value = "M2"
$ws.Range("B2").Value = "Ccl24"
$ws.Range("C2").Value = "Ccr2"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.8467519999999999
$ws.Range("H2").Value = 2.540256
$ws.Range("I2").Value = 1
$ws.Range("J2").Value = 1
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 9.506851333333334
$ws.Range("N2").Value = 28.520554
$ws.Range("O2").Value = 0.04665929098818478
$ws.Range("P2").Value = 0.04665929098818478
$ws.Range("Q2").Value = 8.049945380202667
$ws.Range("R2").Value = 72.44950842182399
$ws.Range("S2").Value = 0.04665929098818478
$ws.Range("T2").Value = 0.04665929098818478

# Row 3: new row for FAPs
$ws.Range("A3").Value = "M2"
$ws.Range("B3").Value = "Ccl24"
$ws.Range("C3").Value = "Ccr2"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.8467519999999999
$ws.Range("H3").Value = 2.540256
$ws.Range("I3").Value = 1
$ws.Range("J3").Value = 1
$ws.Range("K3").Value = 1
$ws.Range("L3").Value = 0.3333333333333333
$ws.Range("M3").Value = 0.01112833333333333
$ws.Range("N3").Value = 0.033385
$ws.Range("O3").Value = 0.00005461746744612846
$ws.Range("P3").Value = 0.00005461746744612846
$ws.Range("Q3").Value = 0.009422938506666666
$ws.Range("R3").Value = 0.08480644655999998
$ws.Range("S3").Value = 0.00005461746744612846
$ws.Range("T3").Value = 0.00005461746744612846

# Row 4: new row for M2 (target cluster self)
$ws.Range("A4").Value = "M2"
$ws.Range("B4").Value = "Ccl24"
$ws.Range("C4").Value = "Ccr2"
$ws.Range("D4").Value = "M2"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.8467519999999999
$ws.Range("H4").Value = 2.540256
$ws.Range("I4").Value = 1
$ws.Range("J4").Value = 1
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 194.1975953333333
$ws.Range("N4").Value = 582.592786
$ws.Range("O4").Value = 0.9531149475424379
$ws.Range("P4").Value = 0.9531149475424379
$ws.Range("Q4").Value = 164.4372022436906
$ws.Range("R4").Value = 1479.934820193216
$ws.Range("S4").Value = 0.9531149475424379
$ws.Range("T4").Value = 0.9531149475424379

# Row 5: new row for sCs
$ws.Range("A5").Value = "M2"
$ws.Range("B5").Value = "Ccl24"
$ws.Range("C5").Value = "Ccr2"
$ws.Range("D5").Value = "sCs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 0.8467519999999999
$ws.Range("H5").Value = 2.540256
$ws.Range("I5").Value = 1
$ws.Range("J5").Value = 1
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.03487066666666667
$ws.Range("N5").Value = 0.104612
$ws.Range("O5").Value = 0.0001711440019312383
$ws.Range("P5").Value = 0.0001711440019312383
$ws.Range("Q5").Value = 0.02952680674133333
$ws.Range("R5").Value = 0.265741260672
$ws.Range("S5").Value = 0.0001711440019312383
$ws.Range("T5").Value = 0.0001711440019312383
